$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "42.263.62"
$ws.Range("E2").Value = "  +0.35%  "
Set-TextValue "D3" "2.269.00"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "307.11"
$ws.Range("E5").Value = "  +0.68%  "
Set-TextValue "D6" "96.94"
$ws.Range("E6").Value = "  +3.11%  "
Set-TextValue "D7" "0.526"
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue "D9" "0.495"
$ws.Range("E9").Value = "  +1.07%  "
Set-TextValue "D10" "35.33"
$ws.Range("E10").Value = "  +4.01%  "
Set-TextValue "D11" "0.0790"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("E12").Value = "  +0.02%  "
Set-TextValue "D13" "6.85"
$ws.Range("E13").Value = "  +2.85%  "
Set-TextValue "D14" "2.619.01"
$ws.Range("E14").Value = "  -0.53%  "
Set-TextValue "D15" "14.82"
$ws.Range("E15").Value = "  +3.23%  "
Set-TextValue "D16" "2.272.78"
$ws.Range("E16").Value = "  +0.10%  "
Set-TextValue "D17" "0.795"
$ws.Range("E17").Value = "  +0.35%  "
Set-TextValue "D18" "42.115.70"
$ws.Range("E18").Value = "  +0.19%  "
Set-TextValue "D19" "12.42"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  +0.62%  "
Set-TextValue "D22" "68.22"
$ws.Range("E22").Value = "  +0.32%  "
Set-TextValue "D23" "238.39"
$ws.Range("E23").Value = "  -2.18%  "
Set-TextValue "D24" "2.57"
$ws.Range("E24").Value = "  -1.06%  "
Set-TextValue "D25" "1.94"
$ws.Range("E25").Value = "  +0.21%  "
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  +0.02%  "
Set-TextValue "D27" "23.65"
$ws.Range("E27").Value = "  -1.75%  "
Set-TextValue "D28" "37.67"
$ws.Range("E28").Value = "  +5.00%  "
Set-TextValue "D29" "9.50"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("E30").Value = "  +0.94%  "
Set-TextValue "D31" "161.94"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("E33").Value = "  +0.04%  "
Set-TextValue "D34" "3.18"
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("E35").Value = "  -2.11%  "
Set-TextValue "D36" "17.16"
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("E38").Value = "  -3.35%  "
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("E40").Value = "  -1.51%  "
Set-TextValue "D41" "4.04"
$ws.Range("E41").Value = "  -4.19%  "
$ws.Range("E42").Value = "  +2.04%  "
Set-TextValue "D43" "1.946.29"
$ws.Range("E43").Value = "  -3.73%  "
Set-TextValue "D44" "18.94"
$ws.Range("E44").Value = "  -3.89%  "
Set-TextValue "D45" "0.0282"
$ws.Range("E45").Value = "  -0.73%  "
Set-TextValue "D46" "9.96"
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("E47").Value = "  -1.17%  "
Set-TextValue "D48" "53.61"
$ws.Range("E48").Value = "  +0.28%  "
Set-TextValue "D49" "71.91"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  -1.66%  "
